$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("writing")
foreach ($r in 43..46) {
  $addr = "N$r"
  try {
    $ws.Range($addr).FormulaArray = "=1+1"
    Write-Output "$addr OK -> $($ws.Range($addr).Formula)"
  } catch {
    Write-Output "$addr FAILED: $_"
  }
}
